$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.000.14"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "1.829.25"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.04"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6219"
$ws.Range("E6").Value = "  -6.53%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07532"
$ws.Range("E8").Value = "  +0.73%  "
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.40"
$ws.Range("E10").Value = "  -4.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07724"
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").Value = "1.833.34"
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.936"
$ws.Range("E13").Value = "  -1.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6630"
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009994"
$ws.Range("E15").Value = "  +14.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.39"
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.026"
$ws.Range("E17").Value = "  -2.78%  "
$ws.Range("D18").Value = "29.019.39"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "225.35"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.29"
$ws.Range("E20").Value = "  -2.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9993"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.146"
$ws.Range("E22").Value = "  -1.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.439"
$ws.Range("E25").Value = "  -2.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1368"
$ws.Range("E26").Value = "  -2.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.89"
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.493"
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.080"
$ws.Range("E29").Value = "  -1.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.020"
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.196"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05188"
$ws.Range("E32").Value = "  -3.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.844"
$ws.Range("E33").Value = "  -0.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7361"
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("E35").Value = "  -2.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.696"
$ws.Range("E36").Value = "  +1.89%  "
$ws.Range("D37").Value = "1.244.78"
$ws.Range("E37").Value = "  -4.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.761"
$ws.Range("E38").Value = "  +0.18%  "
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.300"
$ws.Range("E40").Value = "  -1.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8944"
$ws.Range("E41").Value = "  -1.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.26"
$ws.Range("E43").Value = "  -2.38%  "
$ws.Range("E44").Value = "  +2.28%  "
$ws.Range("D45").Value = "1.982.51"
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.88"
$ws.Range("E46").Value = "  -2.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5107"
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4014"
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.845"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05744"
$ws.Range("E50").Value = "  -2.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.633"
$ws.Range("E51").Value = "  -7.09%  "
